$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Grocery")
$ws.Range("A4").Value = "Vegetables"
$ws.Activate()
$ws.Range("A4").Select() | Out-Null
